# reimplemented the tests with correct url
# Updates per-run metric values across the "runs", "per_resource_all_runs",
# "summary_by_type" and "co2" sheets to reflect a re-run of the Lighthouse/
# per-resource evidence collection against the corrected target URLs
# (e.g. swapped tracked images/scripts for the right resources).

$wb = $excel.ActiveWorkbook
$ws_runs = $wb.Worksheets.Item("runs")
$ws_runs.Range("B2").Value = 15297.624999999995
$ws_runs.Range("C2").Value = 9379.411
$ws_runs.Range("D2").Value = 72
$ws_runs.Range("E2").Value = 15321.006699999994
$ws_runs.Range("F2").Value = 98
$ws_runs.Range("H2").Value = 9379.411
$ws_runs.Range("I2").Value = 56
$ws_runs.Range("J2").Value = 2297354
$ws_runs.Range("B3").Value = 14265.7335
$ws_runs.Range("C3").Value = 8316.996500000001
$ws_runs.Range("D3").Value = 81
$ws_runs.Range("E3").Value = 14273.51215
$ws_runs.Range("F3").Value = 77
$ws_runs.Range("H3").Value = 8316.996500000001
$ws_runs.Range("J3").Value = 2046789
$ws_runs.Range("B4").Value = 13352.910999999995
$ws_runs.Range("C4").Value = 8158.742999999999
$ws_runs.Range("D4").Value = 124
$ws_runs.Range("E4").Value = 13368.503699999994
$ws_runs.Range("F4").Value = 119
$ws_runs.Range("H4").Value = 8158.742999999999
$ws_runs.Range("I4").Value = 56
$ws_runs.Range("J4").Value = 2046730
$ws_runs.Range("B5").Value = 14265.7335
$ws_runs.Range("C5").Value = 8316.996500000001
$ws_runs.Range("D5").Value = 81
$ws_runs.Range("E5").Value = 14273.51215
$ws_runs.Range("F5").Value = 98
$ws_runs.Range("H5").Value = 8316.996500000001
$ws_runs.Range("I5").Value = 56
$ws_runs.Range("J5").Value = 2046789
$ws_per_resource_all_runs = $wb.Worksheets.Item("per_resource_all_runs")
$ws_per_resource_all_runs.Range("E2").Value = 18990
$ws_per_resource_all_runs.Range("F2").Value = 58295
$ws_per_resource_all_runs.Range("G2").Value = 0.005522538869999999
$ws_per_resource_all_runs.Range("H2").Value = 0.007274715786
$ws_per_resource_all_runs.Range("E3").Value = 28150
$ws_per_resource_all_runs.Range("G3").Value = 0.008186385949999998
$ws_per_resource_all_runs.Range("H3").Value = 0.01078374141
$ws_per_resource_all_runs.Range("E4").Value = 26223
$ws_per_resource_all_runs.Range("G4").Value = 0.007625989298999999
$ws_per_resource_all_runs.Range("H4").Value = 0.010045543552200002
$ws_per_resource_all_runs.Range("E5").Value = 12284
$ws_per_resource_all_runs.Range("G5").Value = 0.0035723468919999996
$ws_per_resource_all_runs.Range("H5").Value = 0.0047057719176
$ws_per_resource_all_runs.Range("E11").Value = 43667
$ws_per_resource_all_runs.Range("G11").Value = 0.012698931270999999
$ws_per_resource_all_runs.Range("H11").Value = 0.0167280154938
$ws_per_resource_all_runs.Range("E12").Value = 5180
$ws_per_resource_all_runs.Range("G12").Value = 0.0015064113399999998
$ws_per_resource_all_runs.Range("H12").Value = 0.0019843616520000005
$ws_per_resource_all_runs.Range("E13").Value = 5347
$ws_per_resource_all_runs.Range("G13").Value = 0.0015549771109999997
$ws_per_resource_all_runs.Range("H13").Value = 0.0020483362458000006
$ws_per_resource_all_runs.Range("E14").Value = 3148
$ws_per_resource_all_runs.Range("G14").Value = 0.0009154793239999999
$ws_per_resource_all_runs.Range("H14").Value = 0.0012059402472000002
$ws_per_resource_all_runs.Range("E17").Value = 3173
$ws_per_resource_all_runs.Range("G17").Value = 0.0009227496489999998
$ws_per_resource_all_runs.Range("H17").Value = 0.0012155172822000003
$ws_per_resource_all_runs.Range("E21").Value = 21145
$ws_per_resource_all_runs.Range("G21").Value = 0.006149240885
$ws_per_resource_all_runs.Range("H21").Value = 0.008100256203
$ws_per_resource_all_runs.Range("B26").Value = "https://www.lorenskog.kommune.no/handlers/bv.ashx/ie2aac308-528b-4c39-bbfb-48170d3a0acd/w1920/h1080/k1a9b5eb19e2a/trygg-i-ungdomsland-banner-til-nettside-beskjart.png"
$ws_per_resource_all_runs.Range("D26").Value = "image/png"
$ws_per_resource_all_runs.Range("E26").Value = 71440
$ws_per_resource_all_runs.Range("F26").Value = 71024
$ws_per_resource_all_runs.Range("G26").Value = 0.02077568072
$ws_per_resource_all_runs.Range("H26").Value = 0.027367335216
$ws_per_resource_all_runs.Range("B27").Value = "https://www.lorenskog.kommune.no/handlers/bv.ashx/i124a37e6-3c89-4de1-9076-7342317e4cfd/w1920/h1080/kaa39f6713f54/egenberedskap.jpg"
$ws_per_resource_all_runs.Range("E27").Value = 56241
$ws_per_resource_all_runs.Range("F27").Value = 55859
$ws_per_resource_all_runs.Range("G27").Value = 0.016355613932999998
$ws_per_resource_all_runs.Range("H27").Value = 0.021544881017399996
$ws_per_resource_all_runs.Range("B28").Value = "https://www.lorenskog.kommune.no/handlers/bv.ashx/i666a81bd-00f9-49ac-a4fd-7b7f98f3ed85/w1920/h1080/q13228/k2ea7bdc32116/tsk-okt-print-a3-eller-a4.jpg"
$ws_per_resource_all_runs.Range("D28").Value = "image/jpeg"
$ws_per_resource_all_runs.Range("E28").Value = 123388
$ws_per_resource_all_runs.Range("F28").Value = 122993
$ws_per_resource_all_runs.Range("G28").Value = 0.035882834444
$ws_per_resource_all_runs.Range("H28").Value = 0.0472676477832
$ws_per_resource_all_runs.Range("E34").Value = 6114
$ws_per_resource_all_runs.Range("G34").Value = 0.0017780306819999998
$ws_per_resource_all_runs.Range("H34").Value = 0.0023421596796000003
$ws_per_resource_all_runs.Range("E35").Value = 2620
$ws_per_resource_all_runs.Range("G35").Value = 0.0007619300599999999
$ws_per_resource_all_runs.Range("H35").Value = 0.001003673268
$ws_per_resource_all_runs.Range("E36").Value = 2310
$ws_per_resource_all_runs.Range("G36").Value = 0.0006717780299999999
$ws_per_resource_all_runs.Range("H36").Value = 0.000884918034
$ws_per_resource_all_runs.Range("E42").Value = 8878
$ws_per_resource_all_runs.Range("F42").Value = 30057
$ws_per_resource_all_runs.Range("G42").Value = 0.002581837814
$ws_per_resource_all_runs.Range("H42").Value = 0.003400996669200001
$ws_per_resource_all_runs.Range("E54").Value = 182654
$ws_per_resource_all_runs.Range("G54").Value = 0.05311815770199999
$ws_per_resource_all_runs.Range("H54").Value = 0.0699713500356
$ws_per_resource_all_runs.Range("B55").Value = "https://speech2.leseweb.dk/rawfiles/extern2.min.js"
$ws_per_resource_all_runs.Range("B56").Value = "https://speech2.leseweb.dk/rawfiles/vfact2.min.js"
$ws_per_resource_all_runs.Range("E56").Value = 13134
$ws_per_resource_all_runs.Range("F56").Value = 42731
$ws_per_resource_all_runs.Range("G56").Value = 0.0038195379419999997
$ws_per_resource_all_runs.Range("H56").Value = 0.005031391107600001
$ws_per_resource_all_runs.Range("E61").Value = 19006
$ws_per_resource_all_runs.Range("F61").Value = 58295
$ws_per_resource_all_runs.Range("G61").Value = 0.005527191877999999
$ws_per_resource_all_runs.Range("H61").Value = 0.007280845088400001
$ws_per_resource_all_runs.Range("B85").Value = "https://www.lorenskog.kommune.no/handlers/bv.ashx/ie2aac308-528b-4c39-bbfb-48170d3a0acd/w1920/h1080/k1a9b5eb19e2a/trygg-i-ungdomsland-banner-til-nettside-beskjart.png"
$ws_per_resource_all_runs.Range("D85").Value = "image/png"
$ws_per_resource_all_runs.Range("E85").Value = 71440
$ws_per_resource_all_runs.Range("F85").Value = 71024
$ws_per_resource_all_runs.Range("G85").Value = 0.02077568072
$ws_per_resource_all_runs.Range("H85").Value = 0.027367335216
$ws_per_resource_all_runs.Range("B86").Value = "https://www.lorenskog.kommune.no/handlers/bv.ashx/i124a37e6-3c89-4de1-9076-7342317e4cfd/w1920/h1080/kaa39f6713f54/egenberedskap.jpg"
$ws_per_resource_all_runs.Range("E86").Value = 56241
$ws_per_resource_all_runs.Range("F86").Value = 55859
$ws_per_resource_all_runs.Range("G86").Value = 0.016355613932999998
$ws_per_resource_all_runs.Range("H86").Value = 0.021544881017399996
$ws_per_resource_all_runs.Range("B87").Value = "https://www.lorenskog.kommune.no/handlers/bv.ashx/i666a81bd-00f9-49ac-a4fd-7b7f98f3ed85/w1920/h1080/q13228/k2ea7bdc32116/tsk-okt-print-a3-eller-a4.jpg"
$ws_per_resource_all_runs.Range("D87").Value = "image/jpeg"
$ws_per_resource_all_runs.Range("E87").Value = 123388
$ws_per_resource_all_runs.Range("F87").Value = 122993
$ws_per_resource_all_runs.Range("G87").Value = 0.035882834444
$ws_per_resource_all_runs.Range("H87").Value = 0.0472676477832
$ws_per_resource_all_runs.Range("E99").Value = 1054
$ws_per_resource_all_runs.Range("G99").Value = 0.0003065169019999999
$ws_per_resource_all_runs.Range("H99").Value = 0.00040376779559999997
$ws_per_resource_all_runs.Range("E100").Value = 590
$ws_per_resource_all_runs.Range("G100").Value = 0.00017157967
$ws_per_resource_all_runs.Range("H100").Value = 0.00022601802599999998
$ws_per_resource_all_runs.Range("E113").Value = 182654
$ws_per_resource_all_runs.Range("G113").Value = 0.05311815770199999
$ws_per_resource_all_runs.Range("H113").Value = 0.0699713500356
$ws_per_resource_all_runs.Range("B114").Value = "https://speech13.leseweb.dk/rawfiles/extern2.min.js"
$ws_per_resource_all_runs.Range("B115").Value = "https://speech13.leseweb.dk/rawfiles/vfact2.min.js"
$ws_per_resource_all_runs.Range("E115").Value = 13135
$ws_per_resource_all_runs.Range("G115").Value = 0.0038198287549999995
$ws_per_resource_all_runs.Range("H115").Value = 0.005031774189000001
$ws_per_resource_all_runs.Range("E120").Value = 18971
$ws_per_resource_all_runs.Range("F120").Value = 58295
$ws_per_resource_all_runs.Range("G120").Value = 0.005517013422999999
$ws_per_resource_all_runs.Range("H120").Value = 0.0072674372394000005
$ws_per_resource_all_runs.Range("B144").Value = "https://www.lorenskog.kommune.no/handlers/bv.ashx/ie2aac308-528b-4c39-bbfb-48170d3a0acd/w1920/h1080/k1a9b5eb19e2a/trygg-i-ungdomsland-banner-til-nettside-beskjart.png"
$ws_per_resource_all_runs.Range("D144").Value = "image/png"
$ws_per_resource_all_runs.Range("E144").Value = 71440
$ws_per_resource_all_runs.Range("F144").Value = 71024
$ws_per_resource_all_runs.Range("G144").Value = 0.02077568072
$ws_per_resource_all_runs.Range("H144").Value = 0.027367335216
$ws_per_resource_all_runs.Range("B145").Value = "https://www.lorenskog.kommune.no/handlers/bv.ashx/i124a37e6-3c89-4de1-9076-7342317e4cfd/w1920/h1080/kaa39f6713f54/egenberedskap.jpg"
$ws_per_resource_all_runs.Range("E145").Value = 56241
$ws_per_resource_all_runs.Range("F145").Value = 55859
$ws_per_resource_all_runs.Range("G145").Value = 0.016355613932999998
$ws_per_resource_all_runs.Range("H145").Value = 0.021544881017399996
$ws_per_resource_all_runs.Range("B146").Value = "https://www.lorenskog.kommune.no/handlers/bv.ashx/i666a81bd-00f9-49ac-a4fd-7b7f98f3ed85/w1920/h1080/q13228/k2ea7bdc32116/tsk-okt-print-a3-eller-a4.jpg"
$ws_per_resource_all_runs.Range("D146").Value = "image/jpeg"
$ws_per_resource_all_runs.Range("E146").Value = 123388
$ws_per_resource_all_runs.Range("F146").Value = 122993
$ws_per_resource_all_runs.Range("G146").Value = 0.035882834444
$ws_per_resource_all_runs.Range("H146").Value = 0.0472676477832
$ws_per_resource_all_runs.Range("E158").Value = 660
$ws_per_resource_all_runs.Range("G158").Value = 0.00019193657999999999
$ws_per_resource_all_runs.Range("H158").Value = 0.000252833724
$ws_per_resource_all_runs.Range("E159").Value = 961
$ws_per_resource_all_runs.Range("G159").Value = 0.00027947129299999996
$ws_per_resource_all_runs.Range("H159").Value = 0.0003681412254
$ws_per_resource_all_runs.Range("E160").Value = 8879
$ws_per_resource_all_runs.Range("G160").Value = 0.0025821286269999996
$ws_per_resource_all_runs.Range("H160").Value = 0.0034013797505999996
$ws_per_resource_all_runs.Range("E172").Value = 182654
$ws_per_resource_all_runs.Range("G172").Value = 0.05311815770199999
$ws_per_resource_all_runs.Range("H172").Value = 0.0699713500356
$ws_per_resource_all_runs.Range("B173").Value = "https://speech7.leseweb.dk/rawfiles/extern2.min.js"
$ws_per_resource_all_runs.Range("B174").Value = "https://speech7.leseweb.dk/rawfiles/vfact2.min.js"
$ws_per_resource_all_runs.Range("E174").Value = 13135
$ws_per_resource_all_runs.Range("G174").Value = 0.0038198287549999995
$ws_per_resource_all_runs.Range("H174").Value = 0.005031774189000001
$ws_summary_by_type = $wb.Worksheets.Item("summary_by_type")
$ws_summary_by_type.Range("C2").Value = 56967
$ws_summary_by_type.Range("D2").Value = 174885
$ws_summary_by_type.Range("C3").Value = 338088
$ws_summary_by_type.Range("C4").Value = 2117329
$ws_summary_by_type.Range("D4").Value = 8430839
$ws_summary_by_type.Range("C5").Value = 2614685
$ws_summary_by_type.Range("D5").Value = 2711154
$ws_summary_by_type.Range("C6").Value = 1197489
$ws_co2 = $wb.Worksheets.Item("co2")
$ws_co2.Range("B2").Value = 2297354
$ws_co2.Range("C2").Value = 0.668100408802
$ws_co2.Range("D2").Value = 0.8800735866156001
$ws_co2.Range("B3").Value = 2046789
$ws_co2.Range("C3").Value = 0.595232849457
$ws_co2.Range("D3").Value = 0.7840867956245998
$ws_co2.Range("B4").Value = 2046730
$ws_co2.Range("C4").Value = 0.5952156914899999
$ws_co2.Range("D4").Value = 0.7840641938220001
$ws_co2.Range("B5").Value = 2046789
$ws_co2.Range("C5").Value = 0.595232849457
$ws_co2.Range("D5").Value = 0.7840867956245998
